$wb = $excel.ActiveWorkbook

$oldGuid = "6ce50563-b79a-43c2-adbc-da15402179c1"
$newGuid = "0b0d0a19-424d-4f80-b08e-943769bf5182"
$oldHash = "37e33e84fb24e6fc6812e0a7a99b6215ee3871bc"
$newHash = "86b0948d5317ac65ed5104b5460c4b1755457d91"

# Hyperlink target addresses are historical permalinks (pinned to a commit
# SHA) and are left untouched by this change -- only the display text/cell
# text is updated, matching the source diff (no .rels changes).
$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8a535e60295fea502b775b82b16d61f55f4c6bf9/e2e"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "$ghBase/$oldGuid.md", [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md")
$wsOverview.Range("G2").Value = "2016-11-08 23:15:27"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "$ghBase/$oldGuid.md", [Type]::Missing, [Type]::Missing, "$newGuid.md")
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-11-08 23:15:13"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "$ghBase/$oldGuid.md", [Type]::Missing, [Type]::Missing, "$newGuid.md")
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-11-08 23:15:27"
